# Generate Report for Archive
#
# - Status text "Ready for handoff" -> "In Translation" everywhere it
#   appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share that string).
# - The Status column(s) narrow to fit the new (shorter) text:
#   Overview columns E & F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the now-shorter Status columns to match the refreshed content.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
